# "body full need to work on quartiles"
#
# RealLength_Hull (column S) was recomputed and is now (re)synced with
# RealLength_MEC (column T) for every measurement row. Copy T -> S for
# every row that actually holds a measurement (rows with no data in T
# are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $tCell = $ws.Cells.Item($r, 20)   # column T = RealLength_MEC(cm)
    $tVal = $tCell.Value2
    if ($tVal -ne $null -and $tVal -ne "") {
        $sCell = $ws.Cells.Item($r, 19)   # column S = RealLength_Hull(cm)
        $sCell.Value = $tVal
    }
}
